function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" "27.245.29"
Set-TextCell $ws "E2" "  +0.29%  "
Set-TextCell $ws "D3" "1.907.28"
Set-TextCell $ws "E3" "  +0.13%  "
Set-TextCell $ws "E4" "  +0.11%  "
Set-TextCell $ws "D5" "307.63"
Set-TextCell $ws "E5" "  +0.48%  "
Set-TextCell $ws "E6" "  +0.10%  "
Set-TextCell $ws "D7" "0.5254"
Set-TextCell $ws "E7" "  +0.34%  "
Set-TextCell $ws "D8" "0.3817"
Set-TextCell $ws "E8" "  +1.26%  "
Set-TextCell $ws "D9" "0.07308"
Set-TextCell $ws "E9" "  +0.84%  "
Set-TextCell $ws "D10" "21.61"
Set-TextCell $ws "E10" "  +1.98%  "
Set-TextCell $ws "D11" "0.9058"
Set-TextCell $ws "E11" "  +0.33%  "
Set-TextCell $ws "D12" "0.08184"
Set-TextCell $ws "E12" "  -4.02%  "
Set-TextCell $ws "D13" "96.39"
Set-TextCell $ws "E13" "  -0.66%  "
Set-TextCell $ws "D14" "5.368"
Set-TextCell $ws "E14" "  +1.41%  "
Set-TextCell $ws "D15" "1.631.99"
Set-TextCell $ws "E15" "  -14.40%  "
Set-TextCell $ws "E16" "  +0.09%  "
Set-TextCell $ws "D17" "0.000008691"
Set-TextCell $ws "E17" "  +0.65%  "
Set-TextCell $ws "E18" "  +1.29%  "
Set-TextCell $ws "E19" "  +0.07%  "
Set-TextCell $ws "D20" "27.273.55"
Set-TextCell $ws "E20" "  +0.26%  "
Set-TextCell $ws "E21" "  +1.12%  "
Set-TextCell $ws "E22" "  +1.93%  "
Set-TextCell $ws "D23" "6.511"
Set-TextCell $ws "E23" "  +1.10%  "
Set-TextCell $ws "D24" "2.350"
Set-TextCell $ws "E24" "  +2.24%  "
Set-TextCell $ws "D25" "149.94"
Set-TextCell $ws "E25" "  +1.86%  "
Set-TextCell $ws "E26" "  -0.07%  "
Set-TextCell $ws "D27" "1.741"
Set-TextCell $ws "E27" "  -0.44%  "
Set-TextCell $ws "D28" "116.86"
Set-TextCell $ws "E28" "  +1.64%  "
Set-TextCell $ws "D29" "4.847"
Set-TextCell $ws "E29" "  +0.67%  "
Set-TextCell $ws "D30" "4.868"
Set-TextCell $ws "E30" "  -1.10%  "
Set-TextCell $ws "D31" "0.09253"
Set-TextCell $ws "E31" "  -0.33%  "
Set-TextCell $ws "D32" "0.8256"
Set-TextCell $ws "E32" "  +2.46%  "
Set-TextCell $ws "E33" "  +0.48%  "
Set-TextCell $ws "D34" "1.228"
Set-TextCell $ws "E34" "  -1.15%  "
Set-TextCell $ws "D35" "2.990"
Set-TextCell $ws "E35" "  +1.11%  "
Set-TextCell $ws "B36" "MXToken"
Set-TextCell $ws "C36" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D36" "3.359"
Set-TextCell $ws "E36" "  -2.59%  "
Set-TextCell $ws "B37" "RenderToken"
Set-TextCell $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D37" "2.728"
Set-TextCell $ws "E37" "  +4.15%  "
Set-TextCell $ws "D38" "0.5761"
Set-TextCell $ws "E38" "  +0.80%  "
Set-TextCell $ws "E39" "  +0.20%  "
Set-TextCell $ws "D40" "1.082"
Set-TextCell $ws "E40" "  +0.90%  "
Set-TextCell $ws "D41" "9.072"
Set-TextCell $ws "D42" "6.595"
Set-TextCell $ws "E42" "  -0.66%  "
Set-TextCell $ws "D43" "117.75"
Set-TextCell $ws "E43" "  +1.26%  "
Set-TextCell $ws "E44" "  +0.29%  "
Set-TextCell $ws "D45" "0.4916"
Set-TextCell $ws "E45" "  +0.91%  "
Set-TextCell $ws "D46" "10.23"
Set-TextCell $ws "E46" "  +0.21%  "
Set-TextCell $ws "E47" "  +0.12%  "
Set-TextCell $ws "D48" "1.643"
Set-TextCell $ws "E48" "  +1.61%  "
Set-TextCell $ws "D49" "38.78"
Set-TextCell $ws "E49" "  +3.19%  "
Set-TextCell $ws "D50" "64.60"
Set-TextCell $ws "E50" "  +0.54%  "
Set-TextCell $ws "D51" "0.06048"
Set-TextCell $ws "E51" "  +1.64%  "

Write-Output "Applied 91 cell updates"
